$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: /home/redes/Desktop/MartinezRincon-0d5ae8-main => --- rwx r-x
$ws.Range("A3").Value = "/home/redes/Desktop/MartinezRincon-0d5ae8-main"
$ws.Range("B3").Value = "---"
$ws.Range("C3").Value = "rwx"
$ws.Range("D3").Value = "r-x"

# Row 4: /home/redes/Downloads => rwx rwx rwx
$ws.Range("A4").Value = "/home/redes/Downloads"
$ws.Range("B4").Value = "rwx"
$ws.Range("C4").Value = "rwx"
$ws.Range("D4").Value = "rwx"

# Row 5: /home/redes/Downloads/Carpeta para probar => r-x r-x r-x
$ws.Range("A5").Value = "/home/redes/Downloads/Carpeta para probar"
$ws.Range("B5").Value = "r-x"
$ws.Range("C5").Value = "r-x"
$ws.Range("D5").Value = "r-x"

# Copy formatting from existing rows onto the new rows so borders/alignment match
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C2").Copy()
$ws.Range("B3:D5").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B7").Select() | Out-Null
